$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The header columns X..AG (10 columns: "No. of Sites Reverted" ... "Municipality
# Classification") are removed. The last header column, AH ("Status as of July 4,
# 2025"), together with its data-validation dropdown cell in row 2, shifts left
# to become column X. Deleting the entire columns X:AG achieves exactly that:
# AH1/AH2 (and the AH2 data validation) slide into X1/X2, and the sheet
# dimension/validation range update automatically.
$ws.Range("X1:AG2").EntireColumn.Delete()

# Give the new header cell (X1, "Status as of July 4, 2025") the same bold font
# and border used by the rest of the header row, but WITHOUT the centered
# alignment the other header cells use.
$ws.Range("W1").Copy()
$ws.Range("X1").PasteSpecial(-4122)
$ws.Range("X1").HorizontalAlignment = 1
$ws.Range("X1").VerticalAlignment = -4107
$excel.CutCopyMode = $false

# Add borders to every cell of the data row (A2:X2), matching the header's
# border so the whole table is boxed.
$ws.Range("A2:X2").Borders.LineStyle = 1

Write-Host "Done."
